# Insert a new data row at row 27 (pushing existing rows 27..117 down to
# 28..118). Excel's Rows.Insert() copies formatting (e.g. the date style)
# from the row above automatically, matching the original sheet's dimension
# growing from A1:T117 to A1:T118.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("27:27").Insert()

# Populate the newly inserted row 27 with the new record.
$ws.Cells.Item(27, 1).Value = 1
$ws.Cells.Item(27, 2).Value = 'Agrícola del Norte S.A. de Arica'
$ws.Cells.Item(27, 3).Value = 'Arica y Parinacota'
$ws.Cells.Item(27, 4).Value = (Get-Date -Year 2022 -Month 8 -Day 31 -Hour 0 -Minute 0 -Second 0)
$ws.Cells.Item(27, 5).Value = 15
$ws.Cells.Item(27, 6).Value = 'Fruta'
$ws.Cells.Item(27, 7).Value = 100102
$ws.Cells.Item(27, 8).Value = 'Cítricos'
$ws.Cells.Item(27, 9).Value = 100102004
$ws.Cells.Item(27, 10).Value = 'Mandarina'
$ws.Cells.Item(27, 11).Value = 'Murcott'
$ws.Cells.Item(27, 12).Value = 'Segunda'
$ws.Cells.Item(27, 13).Value = 350
$ws.Cells.Item(27, 14).Value = 14000
$ws.Cells.Item(27, 15).Value = 15000
$ws.Cells.Item(27, 16).Value = 14500
$ws.Cells.Item(27, 17).Value = '$/caja 20 kilos'
$ws.Cells.Item(27, 18).Value = 'Región de Coquimbo'
$ws.Cells.Item(27, 19).Value = 725
$ws.Cells.Item(27, 20).Value = 20
